$d = $word.ActiveDocument

# 1. Update the ID placeholder text in the first paragraph, and drop the
#    trailing space run that used to follow it.
$d.Content.Find.Execute("**ID__AFFARS_5350_topic_9__ID**", $false, $false, $false, $false, $false, $true, 1, $false, "**ID__AFFARS_SUBPART_5350_104__ID**", 2)

$p1 = $d.Paragraphs(1)
$pEnd = $p1.Range.End
$trailingSpace = $d.Range($pEnd - 2, $pEnd - 1)
if ($trailingSpace.Text -eq " ") {
    $trailingSpace.Delete()
}

# 2. Adjust the paragraph's left indent (120 -> 225 twips == 6pt -> 11.25pt)
#    and add a paragraph border (5pt space on each side, no line).
$p1 = $d.Paragraphs(1)
$p1.Format.LeftIndent = 11.25
$p1.Format.Borders.DistanceFromTop = 5
$p1.Format.Borders.DistanceFromLeft = 5
$p1.Format.Borders.DistanceFromBottom = 5
$p1.Format.Borders.DistanceFromRight = 5

Write-Host "done"
